$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PBL II")
$ws.Range("B2").Value = "TEST"
